$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (two new reporting quarters),
# shifting the existing quarter columns (old D:K) to F:M.
$ws.Columns("D:E").Insert()

# Copy number formats (date format for the header row, thousands-format for data rows)
# from the now-shifted columns onto the two new columns so every row keeps its formatting.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns with their reported values.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1694600
$ws.Range("E8").Value = 1723700
$ws.Range("D9").Value = 856800
$ws.Range("E9").Value = 824000
$ws.Range("D10").Value = 837800
$ws.Range("E10").Value = 899700
$ws.Range("D12").Value = 9300
$ws.Range("E12").Value = 6000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 179100
$ws.Range("E15").Value = 170600
$ws.Range("D17").Value = 1071100
$ws.Range("E17").Value = 1027000
$ws.Range("D18").Value = 623500
$ws.Range("E18").Value = 696700
$ws.Range("D20").Value = -11400
$ws.Range("E20").Value = -3200
$ws.Range("D21").Value = 791200
$ws.Range("E21").Value = 864100
$ws.Range("D22").Value = 70100
$ws.Range("E22").Value = 69200
$ws.Range("D23").Value = 542000
$ws.Range("E23").Value = 624300
$ws.Range("D24").Value = 219000
$ws.Range("E24").Value = 257900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 323000
$ws.Range("E26").Value = 366400
$ws.Range("D27").Value = 324100
$ws.Range("E27").Value = 369400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -30900
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 11400
$ws.Range("E32").Value = 3200
$ws.Range("D33").Value = 293200
$ws.Range("E33").Value = 369400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 293200
$ws.Range("E35").Value = 369400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 844600
$ws.Range("E41").Value = 967300
$ws.Range("D42").Value = 213800
$ws.Range("E42").Value = 236500
$ws.Range("D43").Value = 972600
$ws.Range("E43").Value = 943000
$ws.Range("D44").Value = 1032700
$ws.Range("E44").Value = 1007600
$ws.Range("D45").Value = 116300
$ws.Range("E45").Value = 118800
$ws.Range("D46").Value = 3180000
$ws.Range("E46").Value = 3273200
$ws.Range("D47").Value = 103600
$ws.Range("E47").Value = 103600
$ws.Range("D48").Value = 9403800
$ws.Range("E48").Value = 9311800
$ws.Range("D49").Value = 147700
$ws.Range("E49").Value = 150800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1649700
$ws.Range("E52").Value = 1477300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 14484800
$ws.Range("E54").Value = 14316700
$ws.Range("D57").Value = 598100
$ws.Range("E57").Value = 609000
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 618200
$ws.Range("E59").Value = 598900
$ws.Range("D60").Value = 1216300
$ws.Range("E60").Value = 1207900
$ws.Range("D61").Value = 5960100
$ws.Range("E61").Value = 5959300
$ws.Range("D62").Value = 695600
$ws.Range("E62").Value = 519000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 7917400
$ws.Range("E66").Value = 7730600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 6186900
$ws.Range("E72").Value = 6203000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 6567400
$ws.Range("E76").Value = 6586100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 293200
$ws.Range("E81").Value = 369400
$ws.Range("D83").Value = 179100
$ws.Range("E83").Value = 170600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 446200
$ws.Range("E89").Value = 732800
$ws.Range("D91").Value = -289600
$ws.Range("E91").Value = -282200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -266400
$ws.Range("E94").Value = -467900
$ws.Range("D96").Value = -309300
$ws.Range("E96").Value = -309200
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -309400
$ws.Range("E100").Value = -309400
$ws.Range("D101").Value = 6900
$ws.Range("E101").Value = -20700
$ws.Range("D102").Value = -122700
$ws.Range("E102").Value = -65200
